# Daily attendance processing - 2026-01-03 20:58:43
#
# Normalise the "Recorded By" column (G): whenever the attendance record
# was touched by dnasr281@gmail.com together with exactly one other
# recorder (e.g. the automated "System" actor, or a shared admin
# account), list dnasr281@gmail.com first so the human editor is always
# the lead entry in the audit trail.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$targetUser = "dnasr281@gmail.com"
$lastRow = $ws.UsedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $text = $cell.Text

    if ($text -eq $null -or $text -eq "") {
        continue
    }

    if (-not $text.Contains($targetUser)) {
        continue
    }

    $parts = $text -split ", "

    if ($parts.Length -eq 2 -and $parts[0] -ne $targetUser) {
        $other = $parts[0]
        if ($parts[1] -eq $targetUser) {
            $cell.Value = $targetUser + ", " + $other
        }
    }
}
